# Work Plan version 2.3 - Divided Task 5.5.2 into four sub-tasks to
# accommodate English, isiZulu, and Kinyarwanda languages.
#
# Row 50 ("Task 5.5.2.2 Integrated Text to Speech Conversion") is renumbered
# to "Task 5.5.2.4" and its effort estimate is revised (6 days @ 0.25 PM
# becomes 3 days @ 0.333 PM). All downstream totals recalculate
# automatically via the existing worksheet formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort by Task")
$ws.Activate()

# Rename task 5.5.2.2 -> 5.5.2.4 (the text carries a leading BOM character,
# as in the original string).
$bom = [char]0xFEFF
$ws.Range("A50").Value = "$bom" + "Task 5.5.2.4 Integrated Text to Speech Conversion"

# Update the duration (days) and the person-months rate for that task.
$ws.Range("B50").Value = 3
$ws.Range("C50").Value = 0.333

# Restore the selection / scroll position used in the saved workbook.
$ws.Range("C51").Select()
$excel.ActiveWindow.ScrollRow = 44
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
